# Bulk-app translation xlsx: split icon_filepath/audio_filepath columns
# into per-language (en/fra) columns on the "Modules_and_forms" sheet,
# and update the saved UI selection state (active sheet/cell) to match
# the authored workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert two new columns (H:I) on sheet1, pushing the old H/I (audio_filepath,
# unique_id) two slots to the right (J:K).
$ws1.Columns("H:I").Insert()

# Rename the old "icon_filepath" (now G1) / "audio_filepath" (now J1) headers
# to be language-specific, and fill in the newly inserted localized columns.
$ws1.Range("G1").Value = "icon_filepath_en"
$ws1.Range("J1").Value = "audio_filepath_en"
$ws1.Range("H1").Value = "icon_filepath_fra"
$ws1.Range("I1").Value = "audio_filepath_fra"

# The two new localized-fra header cells pick up an explicit black font
# color in the authored file.
$ws1.Range("H1:I1").Font.Color = 0

# Restore the saved selection state: sheet2 was previously the active tab
# with D5 selected; afterwards sheet1 is active (selection H1:I1) and
# sheet2's lingering selection moves to D40.
$ws2.Range("D40").Select()
$ws1.Activate()
$ws1.Range("H1:I1").Select()
